$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41, shifting the existing rows 41-52 down to 42-53.
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point.
$ws.Range("A41").Value = 6
$ws.Range("B41").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44504
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 100114007
$ws.Range("G41").Value = "Jengibre"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 230
$ws.Range("K41").Value = 13000
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 13870
$ws.Range("N41").Value = "$/caja 13 kilos"
$ws.Range("O41").Value = "Perú"
$ws.Range("P41").Value = 1067
$ws.Range("Q41").Value = 13
$ws.Range("R41").Value = "Hortaliza"
